$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RF (column I) values for rows 30 through 64 to the new recalculated value
$ws.Range("I30:I64").Value = 2885.909574468085
